$wb = $excel.ActiveWorkbook

# zh-cn sheet: row for "2e96fad2-f836-4936-94d3-270934021bb9" file
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-17 14:35:49"
$wsZhCn.Range("H3").Value = "2016-03-17 14:36:15"

# de-de sheet: row for "2e96fad2-f836-4936-94d3-270934021bb9" file
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-17 14:35:53"
$wsDeDe.Range("H3").Value = "2016-03-17 14:36:21"
